$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1821192052980132
$ws.Range("C2").Value = 0.5827814569536424
$ws.Range("J2").Value = 0.006622516556291391
$ws.Range("P2").Value = 0.1291390728476821
$ws.Range("S2").Value = 0.09933774834437085
$ws.Range("B3").Value = 0.01104972375690608
$ws.Range("C3").Value = 0.02209944751381215
$ws.Range("J3").Value = 0.03314917127071823
$ws.Range("P3").Value = 0.7624309392265194
$ws.Range("S3").Value = 0.1712707182320442
$ws.Range("J4").Value = 0.02439024390243903
$ws.Range("P4").Value = 0.7804878048780488
$ws.Range("S4").Value = 0.1951219512195122
$ws.Range("B6").Value = 0.06103286384976526
$ws.Range("D6").Value = 0.009389671361502348
$ws.Range("F6").Value = 0.03755868544600939
$ws.Range("J6").Value = 0.2582159624413146
$ws.Range("O6").Value = 0.02816901408450704
$ws.Range("Q6").Value = 0.1643192488262911
$ws.Range("R6").Value = 0.04694835680751173
$ws.Range("S6").Value = 0.3943661971830986
$ws.Range("B7").Value = 0.1016949152542373
$ws.Range("D7").Value = 0.01694915254237288
$ws.Range("F7").Value = 0.05084745762711865
$ws.Range("J7").Value = 0.1525423728813559
$ws.Range("O7").Value = 0.03389830508474576
$ws.Range("Q7").Value = 0.1468926553672316
$ws.Range("R7").Value = 0.0847457627118644
$ws.Range("S7").Value = 0.4124293785310734
$ws.Range("B8").Value = 0.09649122807017543
$ws.Range("D8").Value = 0.02850877192982456
$ws.Range("E8").Value = 0.002192982456140351
$ws.Range("F8").Value = 0.06578947368421052
$ws.Range("J8").Value = 0.09868421052631579
$ws.Range("O8").Value = 0.008771929824561403
$ws.Range("Q8").Value = 0.1995614035087719
$ws.Range("R8").Value = 0.08771929824561403
$ws.Range("S8").Value = 0.412280701754386
$ws.Range("B9").Value = 0.10727969348659
$ws.Range("D9").Value = 0.01149425287356322
$ws.Range("F9").Value = 0.06896551724137931
$ws.Range("J9").Value = 0.08812260536398467
$ws.Range("O9").Value = 0.01149425287356322
$ws.Range("Q9").Value = 0.2452107279693486
$ws.Range("R9").Value = 0.09195402298850575
$ws.Range("S9").Value = 0.3754789272030651
$ws.Range("B10").Value = 0.1014598540145985
$ws.Range("D10").Value = 0.01532846715328467
$ws.Range("E10").Value = 0.00145985401459854
$ws.Range("F10").Value = 0.06569343065693431
$ws.Range("J10").Value = 0.1160583941605839
$ws.Range("O10").Value = 0.00948905109489051
$ws.Range("Q10").Value = 0.210948905109489
$ws.Range("R10").Value = 0.09854014598540146
$ws.Range("S10").Value = 0.381021897810219
$ws.Range("F11").Value = 0.003472222222222222
$ws.Range("G11").Value = 0.1701388888888889
$ws.Range("J11").Value = 0.07291666666666667
$ws.Range("K11").Value = 0.2152777777777778
$ws.Range("L11").Value = 0.5243055555555556
$ws.Range("S11").Value = 0.01388888888888889
$ws.Range("G12").Value = 0.7070063694267515
$ws.Range("J12").Value = 0.2165605095541401
$ws.Range("K12").Value = 0.01273885350318471
$ws.Range("L12").Value = 0.04458598726114649
$ws.Range("S12").Value = 0.01910828025477707
$ws.Range("G13").Value = 0.55
$ws.Range("J13").Value = 0.375
$ws.Range("S13").Value = 0.075
$ws.Range("F15").Value = 0.01345291479820628
$ws.Range("H15").Value = 0.1434977578475336
$ws.Range("I15").Value = 0.07623318385650224
$ws.Range("J15").Value = 0.42152466367713
$ws.Range("K15").Value = 0.03139013452914798
$ws.Range("M15").Value = 0.01345291479820628
$ws.Range("O15").Value = 0.04035874439461883
$ws.Range("S15").Value = 0.2600896860986547
$ws.Range("H16").Value = 0.155
$ws.Range("I16").Value = 0.08
$ws.Range("J16").Value = 0.385
$ws.Range("K16").Value = 0.135
$ws.Range("M16").Value = 0.025
$ws.Range("N16").Value = 0.005
$ws.Range("O16").Value = 0.075
$ws.Range("S16").Value = 0.14
$ws.Range("F17").Value = 0.014
$ws.Range("H17").Value = 0.162
$ws.Range("I17").Value = 0.114
$ws.Range("J17").Value = 0.42
$ws.Range("K17").Value = 0.08400000000000001
$ws.Range("M17").Value = 0.016
$ws.Range("O17").Value = 0.052
$ws.Range("S17").Value = 0.138
$ws.Range("F18").Value = 0.01327433628318584
$ws.Range("H18").Value = 0.1327433628318584
$ws.Range("I18").Value = 0.1238938053097345
$ws.Range("J18").Value = 0.4823008849557522
$ws.Range("K18").Value = 0.07079646017699115
$ws.Range("M18").Value = 0.02654867256637168
$ws.Range("O18").Value = 0.04867256637168142
$ws.Range("S18").Value = 0.1017699115044248
$ws.Range("F19").Value = 0.01231884057971015
$ws.Range("H19").Value = 0.2065217391304348
$ws.Range("I19").Value = 0.1021739130434783
$ws.Range("J19").Value = 0.3731884057971014
$ws.Range("K19").Value = 0.0963768115942029
$ws.Range("M19").Value = 0.01376811594202899
$ws.Range("O19").Value = 0.07391304347826087
$ws.Range("S19").Value = 0.1217391304347826
